$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new data point (2025-03-31, 3.5) was added to the top of the series,
# pushing the rest of the table down by one row. Insert a fresh row at
# position 3 so it inherits the number formats of the existing data rows
# (row 2 at that moment), shift the old row-2 data into it, then overwrite
# row 2 with the new observation.
$ws.Rows.Item(3).Insert()

$prevDate = $ws.Range("A2").Value2
$prevValue = $ws.Range("B2").Value2
$ws.Range("A3").Value = $prevDate
$ws.Range("B3").Value = $prevValue

$ws.Range("A2").Value = 45747
$ws.Range("B2").Value = 3.5

# Match the saved selection left behind in the workbook.
$ws.Range("B3").Select() | Out-Null
